$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 changes from the text "R40" to the text "1" (still a text/string
# cell, not a number) — mirrors the new shared-string entry <t>1</t> added
# to xl/sharedStrings.xml and B11's <v> switching to point at it.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
